# Apply the "cryptos list" refresh captured in the commit diff.
# Every changed cell is an inline string in the source workbook, so for
# any replacement text that looks numeric we force the cell to Text format
# first (NumberFormat "@") to stop Excel from silently converting the
# value to a number and dropping things like trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.213.61"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "2.650.89"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.21"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.90"
$ws.Range("E6").Value = "  +3.98%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.386"
$ws.Range("E10").Value = "  +6.75%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.61"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").Value = "3.122.80"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").Value = "64.085.18"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").Value = "2.648.87"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.91"
$ws.Range("E18").Value = "  +4.29%  "
$ws.Range("E19").Value = "  +4.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.77"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.91"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.58"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.37"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("E25").Value = "  +8.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.73"
$ws.Range("E26").Value = "  +5.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.40"
$ws.Range("E27").Value = "  +8.67%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.22"
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "555.74"
$ws.Range("E29").Value = "  +2.33%  "
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.06"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").Value = "0.0₃0853"
$ws.Range("E33").Value = "  +5.97%  "
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("E35").Value = "  +2.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "169.77"
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("E39").Value = "  +4.67%  "
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "165.46"
$ws.Range("E42").Value = "  -7.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.29"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.84"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.13"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0571"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.629"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.01"
$ws.Range("E48").Value = "  +15.26%  "
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0962"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.97"
$ws.Range("E51").Value = "  +1.33%  "
